$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force Excel to store the value as literal text (matching the
    # original inlineStr cell type) instead of auto-coercing numeric-
    # looking strings ("577.95", "173.30", ...) into floating point
    # numbers, which would silently rewrite/round the digits.
    $ws.Range($range).Value = "'" + $text
}

# Row 2 - Bitcoin
Set-TextValue "D2" "67.012.85"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.126.34"
$ws.Range("E3").Value = "  +1.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "577.95"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6 - Solana
Set-TextValue "D6" "173.30"
$ws.Range("E6").Value = "  +2.83%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.44%  "

# Row 9 - Toncoin
Set-TextValue "D9" "6.46"
$ws.Range("E9").Value = "  -2.53%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.51%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +0.11%  "

# Row 12 - ShibaInu
Set-TextValue "D12" "0.0000248"
$ws.Range("E12").Value = "  -0.67%  "

# Row 13 - Avalanche
Set-TextValue "D13" "37.24"
$ws.Range("E13").Value = "  +2.21%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -1.27%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.643.32"
$ws.Range("E15").Value = "  +1.22%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "66.982.68"
$ws.Range("E16").Value = "  +0.25%  "

# Row 17 - Polkadot
Set-TextValue "D17" "7.16"
$ws.Range("E17").Value = "  -0.33%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.125.08"
$ws.Range("E18").Value = "  +1.22%  "

# Row 19 - Chainlink
Set-TextValue "D19" "16.32"
$ws.Range("E19").Value = "  +0.94%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "477.44"
$ws.Range("E20").Value = "  +2.32%  "

# Row 21 - Polygon
Set-TextValue "D21" "0.712"
$ws.Range("E21").Value = "  -0.37%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +5.48%  "

# Row 23 - Litecoin
Set-TextValue "D23" "84.02"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "13.33"
$ws.Range("E24").Value = "  +1.91%  "

# Row 25 - Fetch.AI
Set-TextValue "D25" "2.30"
$ws.Range("E25").Value = "  -2.18%  "

# Row 26 - RenderToken
Set-TextValue "D26" "10.13"
$ws.Range("E26").Value = "  +0.07%  "

# Row 28 - NEARProtocol
Set-TextValue "D28" "7.93"
$ws.Range("E28").Value = "  -0.78%  "

# Row 29 - ImmutableX
Set-TextValue "D29" "2.38"
$ws.Range("E29").Value = "  -0.76%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.53%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "28.63"
$ws.Range("E31").Value = "  +1.35%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  +0.82%  "

# Row 33 - PEPE (contains U+2083 SUBSCRIPT THREE; built via
# [string]::Concat rather than "+" because this engine's "+" operator
# arithmetically adds operands that both look numeric, e.g.
# "0.0" + [char]0x2083 would compute 0.0 + 8323 instead of concatenating)
$pepePrice = [string]::Concat("0.0", [char]0x2083, "0957")
Set-TextValue "D33" $pepePrice
$ws.Range("E33").Value = "  -7.07%  "

# Row 35 - Filecoin
Set-TextValue "D35" "5.87"
$ws.Range("E35").Value = "  -0.33%  "

# Row 36 - Mantle
Set-TextValue "D36" "0.978"
$ws.Range("E36").Value = "  -2.53%  "

# Row 37 - Arweave
Set-TextValue "D37" "47.30"
$ws.Range("E37").Value = "  +0.57%  "

# Row 38 - Stacks
Set-TextValue "D38" "2.07"
$ws.Range("E38").Value = "  -1.86%  "

# Row 39 - OKB
Set-TextValue "D39" "50.23"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40 - TheGraph
$ws.Range("E40").Value = "  -1.49%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +1.20%  "

# Row 42 - Cosmos
Set-TextValue "D42" "8.68"
$ws.Range("E42").Value = "  -0.01%  "

# Row 43 - Maker
Set-TextValue "D43" "2.814.43"
$ws.Range("E43").Value = "  +1.46%  "

# Row 44 - Bittensor
Set-TextValue "D44" "384.24"
$ws.Range("E44").Value = "  +0.48%  "

# Row 45 - VeChain
Set-TextValue "D45" "0.0355"
$ws.Range("E45").Value = "  -1.35%  "

# Row 46 - dogwifhat
$ws.Range("E46").Value = "  -9.26%  "

# Row 47 - Monero
Set-TextValue "D47" "135.73"
$ws.Range("E47").Value = "  +0.45%  "

# Row 48 - USDe
$ws.Range("E48").Value = "  -0.02%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "25.00"
$ws.Range("E49").Value = "  +0.77%  "

# Row 50 - ThetaToken
$ws.Range("E50").Value = "  -1.13%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -0.51%  "
